$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for the rows that changed upon repulling data
$ws.Range("F6").Value = 5
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = -4
$ws.Range("F11").Value = -5
$ws.Range("F17").Value = -3
$ws.Range("F20").Value = 1
$ws.Range("F21").Value = -5
$ws.Range("F22").Value = 4
